$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.593.93'
$ws.Range('E2').Value = '  -2.30%  '
$ws.Range('D3').Value = '2.895.40'
$ws.Range('E3').Value = '  -2.05%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.39'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.80%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('D9').Value = '2.893.19'
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.99'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.92%  '
$ws.Range('E11').Value = '  -2.55%  '
$ws.Range('E12').Value = '  -2.30%  '
$ws.Range('E13').Value = '  -0.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('D16').Value = '3.374.85'
$ws.Range('E16').Value = '  -2.06%  '
$ws.Range('D17').Value = '61.581.66'
$ws.Range('E17').Value = '  -2.26%  '
$ws.Range('E18').Value = '  -2.06%  '
$ws.Range('D19').Value = '2.885.11'
$ws.Range('E19').Value = '  -2.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '432.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.658'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.84'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.09'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.03'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.34%  '
$ws.Range('E29').Value = '  +3.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.03'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.23%  '
$ws.Range('E31').Value = '  -4.09%  '
$ws.Range('E32').Value = '  -6.15%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.107'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.98%  '
$ws.Range('E35').Value = '  -3.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.960'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.29%  '
$ws.Range('E37').Value = '  -3.38%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.94'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.08%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.83'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -10.12%  '
$ws.Range('E41').Value = '  -1.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.99%  '
$ws.Range('E44').Value = '  -5.01%  '
$ws.Range('D45').Value = '2.710.08'
$ws.Range('E45').Value = '  +0.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '133.02'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '346.33'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.18%  '
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.64'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.94%  '
